$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 397, pushing existing rows 397-439 down to 398-440.
$ws.Rows.Item(397).Insert()

# Populate the newly inserted row 397 with the new data record.
$ws.Range("A397").Value = 10
$ws.Range("B397").Value = 'Vega Modelo de Temuco'
$ws.Range("C397").Value = 'La Araucanía'
$ws.Range("D397").Value = 45223
$ws.Range("E397").Value = 9
$ws.Range("F397").Value = 'Fruta'
$ws.Range("G397").Value = 100103
$ws.Range("H397").Value = 'Frutos de hueso (carozo)'
$ws.Range("I397").Value = 100103004
$ws.Range("J397").Value = 'Durazno'
$ws.Range("K397").Value = 'Florida King'
$ws.Range("L397").Value = 'Primera'
$ws.Range("M397").Value = 45
$ws.Range("N397").Value = 24000
$ws.Range("O397").Value = 24000
$ws.Range("P397").Value = 24000
$ws.Range("Q397").Value = '$/bandeja 10 kilos granel'
$ws.Range("R397").Value = 'Provincia de Limarí'
$ws.Range("S397").Value = 2400
$ws.Range("T397").Value = 10

# Ensure the D397 cell keeps the date/time number format used by the rest of column D.
$ws.Range("D397").NumberFormat = $ws.Range("D398").NumberFormat
